# 7/14 update: when opening the DatePicker, there is not going to be any
# date selected if the user did not pick one before.
#
# - Drop the obsolete "setIsVisible" attribute row.
# - Add "onCancel" (the action fired when pressing cancel).
# - Expand the "onConfirm" description to explain the single/range argument
#   shapes, wrapped onto several lines.
# - Add "startDate" / "endDate" attributes (first/last date initially
#   selected when the picker opens).
# - Tidy header casing ("Value" / "Discription").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse the old trailing blank row 13 into row 14 (which carries the
# bottom-border / thick-bottom formatting) so the table again ends in a
# single, properly formatted last row.
$ws.Rows.Item(13).Delete()

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "Attributes"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Discription"

# --- Attribute rows --------------------------------------------------------
$ws.Range("A2").Value = "isVisible"
$ws.Range("B2").Value = "boolean"
$ws.Range("C2").Value = "是否開啟視窗"

$ws.Range("A3").Value = "mode"
$ws.Range("B3").Value = "string"
$ws.Range("C3").Value = "模式選擇，單日:'single' 多日: 'range'"

$ws.Range("A4").Value = "onCancel"
$ws.Range("B4").Value = "function"
$ws.Range("C4").Value = "按下取消要執行的動作"

$ws.Range("A5").Value = "onConfirm"
$ws.Range("B5").Value = "function"
$ws.Range("C5").Value = "按下確定時要執行的動作 `nsingle: 傳入 date argument`nrange: 傳入 startDate, endDate 兩個 argument"
$ws.Range("C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 59.25

$ws.Range("A6").Value = "displayDate"
$ws.Range("B6").Value = "Date"
$ws.Range("C6").Value = "開啟視窗時，顯示這個日期所在的月份"

$ws.Range("A7").Value = "minDate"
$ws.Range("B7").Value = "Date"
$ws.Range("C7").Value = "最小的可選日期"

$ws.Range("A8").Value = "maxDate"
$ws.Range("B8").Value = "Date"
$ws.Range("C8").Value = "最大的可選日期"

$ws.Range("A9").Value = "startDate"
$ws.Range("B9").Value = "Date"
$ws.Range("C9").Value = "首次打開畫面時所選取的第一個日期"

$ws.Range("A10").Value = "endDate"
$ws.Range("B10").Value = "Date"
$ws.Range("C10").Value = "首次打開畫面時所選取的最後一個日期"

# Match the author's final selection.
$ws.Range("B11").Select()
